$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.329.49'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '3.261.78'
$ws.Range("E3").Value = '  +2.62%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.77'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.52'
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.262.41'
$ws.Range("E8").Value = '  +2.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.80'
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("E12").Value = '  -3.88%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.03'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").Value = '3.791.06'
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("D16").Value = '66.388.19'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.44'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '3.258.31'
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '505.23'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.755'
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.11'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.65'
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.13'
$ws.Range("E25").Value = '  +2.99%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.21'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.129'
$ws.Range("E30").Value = '  +46.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.04'
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("E32").Value = '  -4.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.98'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.38'
$ws.Range("E37").Value = '  +18.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.73'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").Value = '0.0₃0787'
$ws.Range("E39").Value = '  +15.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '495.35'
$ws.Range("E40").Value = '  -2.17%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.84'
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("E44").Value = '  +3.24%  '
$ws.Range("E45").Value = '  -1.31%  '
$ws.Range("D46").Value = '3.001.43'
$ws.Range("E46").Value = '  +6.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.93'
$ws.Range("E47").Value = '  +3.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.51'
$ws.Range("E48").Value = '  +5.55%  '
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.53'
$ws.Range("E51").Value = '  -3.34%  '
